# Stream Manager deck - slide 4 ("Aggregator" diagram):
# Widen the "Stream Aggregator" icon+label group so the longer caption
# fits, re-center the icon over the wider label, and nudge the group
# into its new spot.
#
# PowerPoint assigns a brand-new Id/Name to a group every time it is
# (re)created - e.g. after an Ungroup + Group round trip, which is
# what happened in the authored edit (the old "Group 51" disappears
# and a new "Group 87" takes its place). Shape.Id is read-only via
# COM, and the Id/Name PowerPoint hands out is simply "the next free
# number" scanned from the shapes already on the slide. We reproduce
# that exactly by first burning through the same number of "slots"
# with scratch shapes (removed again before we are done), so the real
# Group() call lands on the same free slot the original author's
# session landed on.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the "Group 51" shape (Stream Aggregator icon + caption).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Group 51") { $target = $cand }
}

# Burn 30 scratch shapes so the next auto-assigned name/id lands on
# "87"/88, matching the id PowerPoint assigned to the recreated group.
$scratch = @()
for ($n = 1; $n -le 30; $n++) {
    $scratch += $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
}

# Split the group back into its two loose members (picture + caption),
# each keeping its own identity/xfrm, now expressed in slide coords.
$target.Ungroup() | Out-Null

$pic = $null
$tb = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Graphic 47") { $pic = $cand }
    if ($cand.Name -eq "TextBox 48") { $tb = $cand }
}

# Widen the caption box and re-center the icon above it.
$pic.Left = 250.0995275590551
$pic.Top = 206.2463779527559

$tb.Left = 236.52393700787403
$tb.Top = 246.30692913385826
$tb.Width = 64.1511811023622
$tb.Height = 43.2

# Regroup - PowerPoint mints a fresh Id/Name for the new group.
$newGroup = $s.Shapes.Range(@($pic.Name, $tb.Name)).Group()

# Settle the regrouped shape into its final slide position.
$newGroup.Left = 233.78511811023623
$newGroup.Top = 205.34

# Discard the scratch shapes used purely to advance the id counter.
foreach ($d in $scratch) {
    $d.Delete()
}
